$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 21641
$ws.Range("J12").Value = 28662.666
$ws.Range("L12").Value = 28662.666
$ws.Range("N12").Value = -29002.666
$ws.Range("H15").Value = 1225.5518
$ws.Range("I15").Value = 1225.5518
$ws.Range("K15").Value = 3676.6554
$ws.Range("M15").Value = -3507.6554
$ws.Range("H17").Value = 4017.5813
$ws.Range("J17").Value = 4017.5813
$ws.Range("L17").Value = 12052.7439
$ws.Range("N17").Value = -12388.7439
$ws.Range("H40").Value = 4474.028
$ws.Range("J40").Value = 4343.3213
$ws.Range("L40").Value = 4343.3213
$ws.Range("N40").Value = -4693.3213
$ws.Range("H96").Value = 521.8
$ws.Range("I96").Value = 402.25
$ws.Range("K96").Value = 1206.75
$ws.Range("M96").Value = 166.25
$ws.Range("H116").Value = 4756.231
$ws.Range("I116").Value = 4055.7144
$ws.Range("J116").Value = 5573.5
$ws.Range("K116").Value = 4055.7144
$ws.Range("L116").Value = 5573.5
$ws.Range("M116").Value = -613.7143999999998
$ws.Range("N116").Value = -12457.5
$ws.Range("H137").Value = 2978.0312
$ws.Range("J137").Value = 3260.9443
$ws.Range("L137").Value = 9782.832900000001
$ws.Range("N137").Value = -14882.8329
$ws.Range("H138").Value = 2773.6704
$ws.Range("I138").Value = 1744.8889
$ws.Range("J138").Value = 3027.3425
$ws.Range("K138").Value = 5234.6667
$ws.Range("L138").Value = 9082.0275
$ws.Range("M138").Value = -94.66669999999976
$ws.Range("N138").Value = -19362.0275

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3251.82
$ws.Range("I32").Value = 1966.0853
$ws.Range("K32").Value = 1966.0853
$ws.Range("M32").Value = -1679.0853
$ws.Range("H74").Value = 88013.125
$ws.Range("I74").Value = 67215.5
$ws.Range("J74").Value = 150406
$ws.Range("K74").Value = 67215.5
$ws.Range("L74").Value = 150406
$ws.Range("M74").Value = -66341.5
$ws.Range("N74").Value = -152154
$ws.Range("H77").Value = 88013.125
$ws.Range("I77").Value = 67215.5
$ws.Range("J77").Value = 150406
$ws.Range("K77").Value = 336077.5
$ws.Range("L77").Value = 752030
$ws.Range("M77").Value = -331709.5
$ws.Range("N77").Value = -760766
$ws.Range("H97").Value = 32057.777
$ws.Range("I97").Value = 30085
$ws.Range("J97").Value = 36003.332
$ws.Range("K97").Value = 30085
$ws.Range("L97").Value = 36003.332
$ws.Range("M97").Value = -29589
$ws.Range("N97").Value = -36995.332
$ws.Range("H126").Value = 5166.3335
$ws.Range("I126").Value = 5166.3335
$ws.Range("K126").Value = 15499.0005
$ws.Range("M126").Value = -13029.0005
$ws.Range("H132").Value = 2132.1482
$ws.Range("I132").Value = 1702.92
$ws.Range("K132").Value = 5108.76
$ws.Range("M132").Value = -2578.76

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4261.25
$ws.Range("I94").Value = 601.0571
$ws.Range("J94").Value = 14115.615
$ws.Range("K94").Value = 601.0571
$ws.Range("L94").Value = 14115.615
$ws.Range("M94").Value = -150.0571
$ws.Range("N94").Value = -15017.615
$ws.Range("H107").Value = 2567.2222
$ws.Range("I107").Value = 2686.4285
$ws.Range("J107").Value = 2150
$ws.Range("K107").Value = 2686.4285
$ws.Range("L107").Value = 2150
$ws.Range("M107").Value = -766.4285
$ws.Range("N107").Value = -5990
$ws.Range("H134").Value = 2773.7317
$ws.Range("I134").Value = 755.62067
$ws.Range("J134").Value = 7650.8335
$ws.Range("K134").Value = 2266.86201
$ws.Range("L134").Value = 22952.5005
$ws.Range("M134").Value = 268.1379900000002
$ws.Range("N134").Value = -28022.5005

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 18999
$ws.Range("J45").Value = 18999
$ws.Range("L45").Value = 18999
$ws.Range("N45").Value = -20185
$ws.Range("H50").Value = 2614.353
$ws.Range("J50").Value = 2614.353
$ws.Range("L50").Value = 2614.353
$ws.Range("N50").Value = -3864.353
$ws.Range("H60").Value = 7446.5
$ws.Range("I60").Value = 7446.5
$ws.Range("K60").Value = 7446.5
$ws.Range("M60").Value = -6935.5
$ws.Range("H105").Value = 2724.2727
$ws.Range("I105").Value = 1578.1666
$ws.Range("K105").Value = 1578.1666
$ws.Range("M105").Value = 168.8334
$ws.Range("H132").Value = 37836.07
$ws.Range("I132").Value = 41871.48
$ws.Range("K132").Value = 125614.44
$ws.Range("M132").Value = -123084.44
$ws.Range("H134").Value = 27143.27
$ws.Range("I134").Value = 33165.9
$ws.Range("J134").Value = 5311.25
$ws.Range("K134").Value = 99497.70000000001
$ws.Range("L134").Value = 15933.75
$ws.Range("M134").Value = -96962.70000000001
$ws.Range("N134").Value = -21003.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 294.60715
$ws.Range("I2").Value = 160.6
$ws.Range("J2").Value = 369.05554
$ws.Range("K2").Value = 963.5999999999999
$ws.Range("L2").Value = 2214.33324
$ws.Range("M2").Value = -850.5999999999999
$ws.Range("N2").Value = -2440.33324
$ws.Range("H37").Value = 41699.332
$ws.Range("J37").Value = 41699.332
$ws.Range("L37").Value = 125097.996
$ws.Range("N37").Value = -125321.996
$ws.Range("H108").Value = 1659.6428
$ws.Range("I108").Value = 404.375
$ws.Range("J108").Value = 3333.3333
$ws.Range("K108").Value = 1213.125
$ws.Range("L108").Value = 9999.999899999999
$ws.Range("M108").Value = 1666.875
$ws.Range("N108").Value = -15759.9999
$ws.Range("H110").Value = 21989.5
$ws.Range("J110").Value = 21989.5
$ws.Range("L110").Value = 65968.5
$ws.Range("N110").Value = -74148.5
$ws.Range("H111").Value = 1536.2
$ws.Range("I111").Value = 93.666664
$ws.Range("K111").Value = 280.999992
$ws.Range("M111").Value = 2786.000008
$ws.Range("H117").Value = 1869.6
$ws.Range("J117").Value = 1866.3334
$ws.Range("L117").Value = 5599.0002
$ws.Range("N117").Value = -12483.0002
$ws.Range("H121").Value = 12346426
$ws.Range("I121").Value = 461.5
$ws.Range("J121").Value = 22223196
$ws.Range("K121").Value = 1384.5
$ws.Range("L121").Value = 66669588
$ws.Range("M121").Value = -74.5
$ws.Range("N121").Value = -66672208
$ws.Range("H131").Value = 13023213
$ws.Range("I131").Value = 11905595
$ws.Range("J131").Value = 13336146
$ws.Range("K131").Value = 35716785
$ws.Range("L131").Value = 40008438
$ws.Range("M131").Value = -35711745
$ws.Range("N131").Value = -40018518
$ws.Range("H132").Value = 3025.6667
$ws.Range("I132").Value = 2896.6
$ws.Range("J132").Value = 3075.3076
$ws.Range("K132").Value = 26069.4
$ws.Range("L132").Value = 27677.7684
$ws.Range("M132").Value = -23539.4
$ws.Range("N132").Value = -32737.7684
$ws.Range("H137").Value = 3330.2727
$ws.Range("I137").Value = 1663.1666
$ws.Range("J137").Value = 5330.8
$ws.Range("K137").Value = 4989.4998
$ws.Range("L137").Value = 15992.4
$ws.Range("M137").Value = 110.5002000000004
$ws.Range("N137").Value = -26192.4
$ws.Range("H140").Value = 2335.7
$ws.Range("I140").Value = 1928.6666
$ws.Range("K140").Value = 5785.9998
$ws.Range("M140").Value = -605.9997999999996

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 43278.332
$ws.Range("J118").Value = 43278.332
$ws.Range("L118").Value = 43278.332
$ws.Range("N118").Value = -46592.332
$ws.Range("H126").Value = 3464.6667
$ws.Range("I126").Value = 3311.8572
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 9935.571599999999
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -7465.571599999999
$ws.Range("N126").Value = -16938.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5498.6
$ws.Range("I7").Value = 2664.8333
$ws.Range("K7").Value = 2664.8333
$ws.Range("M7").Value = -2552.8333
$ws.Range("H40").Value = 5825.2354
$ws.Range("I40").Value = 2944.0908
$ws.Range("J40").Value = 11107.333
$ws.Range("K40").Value = 2944.0908
$ws.Range("L40").Value = 11107.333
$ws.Range("M40").Value = -2808.0908
$ws.Range("N40").Value = -11379.333
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9309
$ws.Range("H126").Value = 5498.6
$ws.Range("I126").Value = 2664.8333
$ws.Range("K126").Value = 7994.499899999999
$ws.Range("M126").Value = -5524.499899999999
$ws.Range("H132").Value = 7085.9443
$ws.Range("I132").Value = 7431.8213
$ws.Range("J132").Value = 5875.375
$ws.Range("K132").Value = 22295.4639
$ws.Range("L132").Value = 17626.125
$ws.Range("M132").Value = -19765.4639
$ws.Range("N132").Value = -22686.125
$ws.Range("H136").Value = 44571.832
$ws.Range("J136").Value = 7099.9
$ws.Range("L136").Value = 21299.7
$ws.Range("N136").Value = -26399.7

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1404.5
$ws.Range("I113").Value = 991.25
$ws.Range("K113").Value = 2973.75
$ws.Range("M113").Value = -803.75
